$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '74.595.12'
Set-TextValue 'E2' '  +8.99%  '
Set-TextValue 'D3' '2.592.29'
Set-TextValue 'E3' '  +6.80%  '
Set-TextValue 'E4' '  +0.08%  '
Set-TextValue 'D5' '186.11'
Set-TextValue 'E5' '  +16.23%  '
Set-TextValue 'D6' '580.63'
Set-TextValue 'E6' '  +4.07%  '
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'D8' '0.209'
Set-TextValue 'E8' '  +28.72%  '
Set-TextValue 'D9' '0.535'
Set-TextValue 'E9' '  +5.39%  '
Set-TextValue 'D10' '2.593.22'
Set-TextValue 'E10' '  +6.93%  '
Set-TextValue 'E11' '  -0.25%  '
Set-TextValue 'D12' '0.360'
Set-TextValue 'E12' '  +8.82%  '
Set-TextValue 'D13' '4.81'
Set-TextValue 'E13' '  +3.87%  '
Set-TextValue 'D14' '0.0000193'
Set-TextValue 'E14' '  +11.50%  '
Set-TextValue 'D15' '74.633.93'
Set-TextValue 'E15' '  +9.28%  '
Set-TextValue 'D16' '3.057.93'
Set-TextValue 'E16' '  +6.51%  '
Set-TextValue 'D17' '26.26'
Set-TextValue 'E17' '  +14.09%  '
Set-TextValue 'D18' '2.605.21'
Set-TextValue 'E18' '  +7.40%  '
Set-TextValue 'D19' '8.47'
Set-TextValue 'E19' '  +23.15%  '
Set-TextValue 'D20' '11.71'
Set-TextValue 'E20' '  +12.08%  '
Set-TextValue 'D21' '378.44'
Set-TextValue 'E21' '  +13.03%  '
Set-TextValue 'D22' '2.30'
Set-TextValue 'E22' '  +20.88%  '
Set-TextValue 'D23' '4.07'
Set-TextValue 'E23' '  +6.60%  '
Set-TextValue 'E24' '  +0.00%  '
Set-TextValue 'D25' '70.05'
Set-TextValue 'E25' '  +5.17%  '
Set-TextValue 'D26' '4.20'
Set-TextValue 'E26' '  +14.27%  '
Set-TextValue 'D27' '9.18'
Set-TextValue 'E27' '  +11.96%  '
Set-TextValue 'D28' '2.734.95'
Set-TextValue 'E28' '  +7.24%  '
Set-TextValue 'E29' '  -0.56%  '
Set-TextValue 'D30' '0.0₃0944'
Set-TextValue 'E30' '  +15.84%  '
Set-TextValue 'D31' '7.95'
Set-TextValue 'E31' '  +11.49%  '
Set-TextValue 'D32' '501.43'
Set-TextValue 'E32' '  +17.78%  '
Set-TextValue 'D33' '1.33'
Set-TextValue 'E33' '  +16.45%  '
Set-TextValue 'E34' '  +6.31%  '
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  +0.04%  '
Set-TextValue 'D36' '0.120'
Set-TextValue 'E36' '  +14.39%  '
Set-TextValue 'D37' '159.25'
Set-TextValue 'E37' '  +0.33%  '
Set-TextValue 'D38' '19.23'
Set-TextValue 'E38' '  +7.46%  '
Set-TextValue 'E40' '  +0.00%  '
Set-TextValue 'E41' '  +15.35%  '
Set-TextValue 'D42' '1.68'
Set-TextValue 'E42' '  +13.29%  '
Set-TextValue 'E43' '  +8.02%  '
Set-TextValue 'D44' '2.43'
Set-TextValue 'E44' '  +20.73%  '
Set-TextValue 'D45' '39.10'
Set-TextValue 'E45' '  +4.94%  '
Set-TextValue 'D46' '1.16'
Set-TextValue 'E46' '  +8.13%  '
Set-TextValue 'D47' '148.58'
Set-TextValue 'E47' '  +11.84%  '
Set-TextValue 'D48' '0.0816'
Set-TextValue 'E48' '  +14.43%  '
Set-TextValue 'D49' '3.63'
Set-TextValue 'E49' '  +8.87%  '
Set-TextValue 'D50' '0.520'
Set-TextValue 'E50' '  +8.29%  '
Set-TextValue 'D51' '0.583'
Set-TextValue 'E51' '  +4.82%  '
